$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header text values in the exact order the new vocabulary was
# introduced (B1, C1, A1, D1, E1, F1, G1, H1, I1, J1, K1) so that the
# regenerated shared-strings table lands on the same indices as the
# authored workbook.
$ws.Cells.Item(1, 2).Value = "SeqBaseT"
$ws.Cells.Item(1, 3).Value = "ParBaseT"
$ws.Cells.Item(1, 1).Value = "NumFish"
$ws.Cells.Item(1, 4).Value = "MPISeqBaseT"
$ws.Cells.Item(1, 5).Value = "MPIParBaseT"
$ws.Cells.Item(1, 6).Value = "NumThreads"
$ws.Cells.Item(1, 7).Value = "ParBaseT"
$ws.Cells.Item(1, 8).Value = "MPIParBaseT"
$ws.Cells.Item(1, 9).Value = "NumNodes"
$ws.Cells.Item(1, 10).Value = "MPISeqBaseT"
$ws.Cells.Item(1, 11).Value = "MPIParBaseT"

# Data rows 2-5, columns A-K (row 5 / cols I,J,K intentionally left blank)
$data = @(
    @(10000,    3.2684419999999998, 3.3780760000000001, 0.96280500000000002, 1.0972649999999999, 2,  331.269969,         89.734123999999994, 2,     175.07313300000001, 172.31019599999999),
    @(100000,   32.742747999999999, 33.122490999999997, 9.0091509999999992,  9.1660839999999997,  4,  327.678245,         89.760159999999999, 3,     116.310136,         117.08775799999999),
    @(1000000,  327.678245,         330.709878,         88.853880000000004,  89.760159999999999,  8,  331.49715800000001, 90.355018000000001, 4,     88.853880000000004, 89.760159999999999),
    @(10000000, 3253.958807,        3317.2062340000002, 958.066688,          965.17265699999996,  16, 332.09230500000001, 90.144876999999994, $null, $null,               $null)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($r + 2, $c + 1).Value = $val
        }
    }
}

# Column widths to match the new layout (values back-solved so the
# engine's width round-trip lands as close as possible to the widths
# Excel's own best-fit produced in the authored workbook).
$ws.Columns.Item(1).ColumnWidth = 8.833333
$ws.Columns.Item(2).ColumnWidth = 9.333333
$ws.Columns.Item(4).ColumnWidth = 12.0
$ws.Columns.Item(5).ColumnWidth = 11.666667
$ws.Columns.Item(6).ColumnWidth = 11.5
$ws.Columns.Item(8).ColumnWidth = 11.666667
$ws.Columns.Item(9).ColumnWidth = 10.166667
$ws.Columns.Item(10).ColumnWidth = 12.0
$ws.Columns.Item(11).ColumnWidth = 11.666667

# Update the selection to match the recorded state
$ws.Range("E16").Select()
